# Fix: sqrt was missing when calculating SE for demo table.
# Replace each "mean ±SE (min to max)" cell value with the corrected SE.

$d = $word.ActiveDocument

$replacements = @(
    @("43.39 ±0.54 (16 to 62)", "43.39 ±2.61 (16 to 62)"),
    @("46.64 ±0.43 (32 to 67)", "46.64 ±2.02 (32 to 67)"),
    @("32.83 ±0.35 (14 to 52)", "32.83 ±1.66 (14 to 52)"),
    @("25.04 ±0.34 (13 to 42)", "25.04 ±1.68 (13 to 42)"),

    @("26.70 ±0.32 (17 to 43)", "26.70 ±1.51 (17 to 43)"),
    @("30.68 ±0.37 (17 to 44)", "30.68 ±1.73 (17 to 44)"),
    @("29.04 ±0.30 (20 to 45)", "29.04 ±1.44 (20 to 45)"),
    @("27.42 ±0.24 (21 to 42)", "27.42 ±1.15 (21 to 42)"),

    @("3.55 ±0.05 (1 to 5)", "3.55 ±0.23 (1 to 5)"),
    @("4.00 ±0.06 (1 to 5)", "4.00 ±0.29 (1 to 5)"),
    @("3.65 ±0.04 (2 to 5)", "3.65 ±0.21 (2 to 5)"),
    @("4.33 ±0.03 (2 to 5)", "4.33 ±0.17 (2 to 5)"),

    @("108.35 ±0.51 (84 to 130)", "108.35 ±2.46 (84 to 130)"),
    @("113.23 ±0.52 (91 to 133)", "113.23 ±2.43 (91 to 133)"),
    @("111.98 ±0.63 (78 to 144)", "111.98 ±3.00 (78 to 144)"),
    @("109.90 ±0.39 (92 to 131)", "109.90 ±1.92 (92 to 131)"),

    @("92.61 ±1.82 (18 to 148)", "92.61 ±8.74 (18 to 148)"),
    @("145.86 ±1.55 (85 to 201)", "145.86 ±7.26 (85 to 201)"),
    @("151.61 ±1.79 (55 to 217)", "151.61 ±8.58 (55 to 217)"),
    @("44.58 ±1.46 (9 to 142)", "44.58 ±7.15 (9 to 142)")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
